# meilleure detection si nom feuille/colonne est conforme

$wb = $excel.ActiveWorkbook

$sheetChem     = $wb.Worksheets.Item("Produits_chimiques")
$sheetMentions = $wb.Worksheets.Item("Mentions_de_danger")
$sheetConseils = $wb.Worksheets.Item("Conseils_de_prudence")
$sheetLink     = $wb.Worksheets.Item("Link")

# --- Produits_chimiques header row (row 1): spaces -> underscores in column-name tokens ---
$sheetChem.Range("A1").Value = "Nom_francais|type='string'"
$sheetChem.Range("B1").Value = "Nom_anglais|type='string'"
$sheetChem.Range("D1").Value = "NCAS|type='string'"
$sheetChem.Range("F1").Value = "Masse_molaire_gmol|type='float'"
$sheetChem.Range("I1").Value = "Mentions_de_danger|reference=Mentions_de_danger"
$sheetChem.Range("J1").Value = "Conseils_de_prudence|reference=Conseils_de_prudence"
$sheetChem.Range("K1").Value = "FDS_piece_jointe"
$sheetChem.Range("L1").Value = "Hazard_statements"
$sheetChem.Range("M1").Value = "Precautionary_statements"
$sheetChem.Range("N1").Value = "MSDS_piece_jointe"
$sheetChem.Range("O1").Value = "Lien_FDS"
$sheetChem.Range("P1").Value = "Lien_MSDS"

# --- Mentions_de_danger header row (row 1) ---
$sheetMentions.Range("B1").Value = "Mentions_de_danger|type='string'"
$sheetMentions.Range("C1").Value = "Hazard_statements"

# --- Conseils_de_prudence header row (row 1) ---
$sheetConseils.Range("B1").Value = "Conseils_de_prudence|type='string'"
$sheetConseils.Range("C1").Value = "Precautionary_statements|type='string'"
$sheetConseils.Rows.Item(1).RowHeight = 17.9

# --- Link sheet data (sheet/column name references) ---
$sheetLink.Range("A2").Value = "Produits_chimiques"
$sheetLink.Range("B2").Value = "Nom_francais"
$sheetLink.Range("C2").Value = "Nom_anglais"
$sheetLink.Range("A3").Value = "Produits_chimiques"
$sheetLink.Range("B3").Value = "Mentions_de_danger"
$sheetLink.Range("C3").Value = "Precautionary_statements"
$sheetLink.Range("A4").Value = "Produits_chimiques"
$sheetLink.Range("B4").Value = "Lien_FDS"
$sheetLink.Range("C4").Value = "Lien_MSDS"
$sheetLink.Rows.Item(1).RowHeight = 14.9
$sheetLink.Rows.Item(2).RowHeight = 14.9
$sheetLink.Rows.Item(3).RowHeight = 14.9

# --- View / selection state per sheet ---

# Produits_chimiques: no longer the active tab; scrolled to column M, cell N2 selected
$sheetChem.Activate()
$excel.ActiveWindow.ScrollColumn = 13
$excel.ActiveWindow.ScrollRow = 1
$sheetChem.Range("N2").Select()

# Conseils_de_prudence: scrolled to column C, cell C2 selected
$sheetConseils.Activate()
$excel.ActiveWindow.ScrollColumn = 3
$excel.ActiveWindow.ScrollRow = 1
$sheetConseils.Range("C2").Select()

# Link: cell A2 selected (unchanged scroll position)
$sheetLink.Activate()
$sheetLink.Range("A2").Select()

# Mentions_de_danger becomes the active/selected tab (workbook activeTab index 1),
# scrolled to column C, cell C2 selected
$sheetMentions.Activate()
$excel.ActiveWindow.ScrollColumn = 3
$excel.ActiveWindow.ScrollRow = 1
$sheetMentions.Range("C2").Select()
